$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p>' + `
  '<w:r/>' + `
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
  '<w:r><w:t>: Read our review of Age of the Gods Medusa &amp; Monsters online slot game and play for free. Learn about the special features, RTP, and winning potential.</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$metaPara.Range.InsertXML($metaXml) | Out-Null

# ---------------------------------------------------------------------
# 2) Near the end of the document: the duplicate bold title paragraph
#    ("Play Age of the Gods Medusa & Monsters Free | Review") is
#    removed, and the italic paragraph that follows it has its text
#    replaced with the DALLE image prompt (formatting/run shape kept).
# ---------------------------------------------------------------------
$oldReadOurReview = "Read our review of Age of the Gods Medusa & Monsters online slot game and play for free. Learn about the special features, RTP, and winning potential."
$newPrompt = "Prompt: Create a Cartoon Style Feature Image for Age of the Gods Medusa & Monsters DALLE, please create a cartoon-style feature image for the online slot game Age of the Gods Medusa & Monsters. The image should feature a happy Maya warrior with glasses, in a scene inspired by the game's Greek Mythology theme. The warrior can be seen holding a sword or other weapon, with Medusa's gaze fixed on him in the background, surrounded by rocks and the sea. Please use bright and vivid colors to make the image eye-catching and exciting, and make sure that the image conveys the adventurous nature of the game while also highlighting its Greek Mythology theme. Thank you!"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq ($oldReadOurReview + "`r")) {
        $textRange = $d.Range($p.Range.Start, $p.Range.End - 1)
        $textRange.Text = $newPrompt
        break
    }
}

$boldTitleText = "Play Age of the Gods Medusa & Monsters Free | Review"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($i -ne 1 -and $p.Range.Text -eq ($boldTitleText + "`r")) {
        $p.Range.Delete()
        break
    }
}
